$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 29250
$ws.Range("J10").Value = 29250
$ws.Range("L10").Value = 29250
$ws.Range("N10").Value = -29836
$ws.Range("H53").Value = 404.7586
$ws.Range("I53").Value = 260
$ws.Range("J53").Value = 539.86664
$ws.Range("K53").Value = 260
$ws.Range("L53").Value = 539.86664
$ws.Range("M53").Value = 377
$ws.Range("N53").Value = -1813.86664
$ws.Range("H113").Value = 3685.55
$ws.Range("I113").Value = 1710.6666
$ws.Range("J113").Value = 4531.9287
$ws.Range("K113").Value = 1710.6666
$ws.Range("L113").Value = 4531.9287
$ws.Range("M113").Value = 1543.3334
$ws.Range("N113").Value = -11039.9287
$ws.Range("H116").Value = 423884.25
$ws.Range("I116").Value = 1253000.6
$ws.Range("J116").Value = 9326.0625
$ws.Range("K116").Value = 1253000.6
$ws.Range("L116").Value = 9326.0625
$ws.Range("M116").Value = -1249558.6
$ws.Range("N116").Value = -16210.0625
$ws.Range("H132").Value = 24880022
$ws.Range("I132").Value = 37408744
$ws.Range("J132").Value = 717486.5
$ws.Range("K132").Value = 112226232
$ws.Range("L132").Value = 2152459.5
$ws.Range("M132").Value = -112223702
$ws.Range("N132").Value = -2157519.5
$ws.Range("H135").Value = 227.42857
$ws.Range("I135").Value = 198.81482
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 1789.33338
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = 745.66662
$ws.Range("N135").Value = -14070
$ws.Range("H137").Value = 2183.7925
$ws.Range("I137").Value = 1028.0256
$ws.Range("J137").Value = 5403.4287
$ws.Range("K137").Value = 3084.0768
$ws.Range("L137").Value = 16210.2861
$ws.Range("M137").Value = -534.0767999999998
$ws.Range("N137").Value = -21310.2861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 18341.334
$ws.Range("I36").Value = 10012
$ws.Range("K36").Value = 10012
$ws.Range("M36").Value = -9666
$ws.Range("H61").Value = 1364.762
$ws.Range("I61").Value = 1066.3158
$ws.Range("J61").Value = 4200
$ws.Range("K61").Value = 1066.3158
$ws.Range("L61").Value = 4200
$ws.Range("M61").Value = -854.3158000000001
$ws.Range("N61").Value = -4624
$ws.Range("H132").Value = 1728.475
$ws.Range("I132").Value = 1336.742
$ws.Range("J132").Value = 3077.7778
$ws.Range("K132").Value = 4010.226
$ws.Range("L132").Value = 9233.3334
$ws.Range("M132").Value = -1480.226
$ws.Range("N132").Value = -14293.3334
$ws.Range("H136").Value = 1364.762
$ws.Range("I136").Value = 1066.3158
$ws.Range("J136").Value = 4200
$ws.Range("K136").Value = 3198.9474
$ws.Range("L136").Value = 12600
$ws.Range("M136").Value = -648.9474
$ws.Range("N136").Value = -17700

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3284.4
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 3730.5
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 3730.5
$ws.Range("M99").Value = -2
$ws.Range("N99").Value = -6726.5
$ws.Range("H134").Value = 1586.25
$ws.Range("I134").Value = 1148.3
$ws.Range("J134").Value = 3776
$ws.Range("K134").Value = 3444.9
$ws.Range("L134").Value = 11328
$ws.Range("M134").Value = -909.8999999999996
$ws.Range("N134").Value = -16398

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9806062
$ws.Range("I31").Value = 1140.6757
$ws.Range("K31").Value = 1140.6757
$ws.Range("M31").Value = -845.6757
$ws.Range("H34").Value = 9806062
$ws.Range("I34").Value = 1140.6757
$ws.Range("K34").Value = 1140.6757
$ws.Range("M34").Value = -938.6757
$ws.Range("H62").Value = 10000
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -9376
$ws.Range("H65").Value = 10000
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 50000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -46880
$ws.Range("H68").Value = 86249.25
$ws.Range("J68").Value = 86249.25
$ws.Range("L68").Value = 86249.25
$ws.Range("N68").Value = -87747.25
$ws.Range("H71").Value = 86249.25
$ws.Range("J71").Value = 86249.25
$ws.Range("L71").Value = 258747.75
$ws.Range("N71").Value = -266235.75
$ws.Range("H99").Value = 10005996
$ws.Range("I99").Value = 20004192
$ws.Range("J99").Value = 7800
$ws.Range("K99").Value = 20004192
$ws.Range("L99").Value = 7800
$ws.Range("M99").Value = -20002694
$ws.Range("N99").Value = -10796
$ws.Range("H126").Value = 10005996
$ws.Range("I126").Value = 20004192
$ws.Range("J126").Value = 7800
$ws.Range("K126").Value = 60012576
$ws.Range("L126").Value = 23400
$ws.Range("M126").Value = -60010106
$ws.Range("N126").Value = -28340
$ws.Range("H132").Value = 1686.0625
$ws.Range("I132").Value = 993.5833
$ws.Range("J132").Value = 3763.5
$ws.Range("K132").Value = 2980.7499
$ws.Range("L132").Value = 11290.5
$ws.Range("M132").Value = -450.7498999999998
$ws.Range("N132").Value = -16350.5
$ws.Range("H134").Value = 2046.4412
$ws.Range("I134").Value = 893.2778
$ws.Range("J134").Value = 3343.75
$ws.Range("K134").Value = 2679.8334
$ws.Range("L134").Value = 10031.25
$ws.Range("M134").Value = -144.8334
$ws.Range("N134").Value = -15101.25
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2358.682
$ws.Range("I132").Value = 1504.9375
$ws.Range("J132").Value = 4635.3335
$ws.Range("K132").Value = 4514.8125
$ws.Range("L132").Value = 13906.0005
$ws.Range("M132").Value = -1984.8125
$ws.Range("N132").Value = -18966.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11051.82
$ws.Range("I132").Value = 15149.889
$ws.Range("J132").Value = 7539.1904
$ws.Range("K132").Value = 45449.667
$ws.Range("L132").Value = 22617.5712
$ws.Range("M132").Value = -42919.667
$ws.Range("N132").Value = -27677.5712

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 629079.5600000001
$ws.Range("I126").Value = 2547.375
$ws.Range("J126").Value = 1185997.1
$ws.Range("K126").Value = 7642.125
$ws.Range("L126").Value = 3557991.3
$ws.Range("M126").Value = -5172.125
$ws.Range("N126").Value = -3562931.3
$ws.Range("H132").Value = 5557133
$ws.Range("I132").Value = 1132.5834
$ws.Range("J132").Value = 27781136
$ws.Range("K132").Value = 3397.7502
$ws.Range("L132").Value = 83343408
$ws.Range("M132").Value = -867.7501999999999
$ws.Range("N132").Value = -83348468
$ws.Range("H136").Value = 2560.1
$ws.Range("I136").Value = 635.3929000000001
$ws.Range("J136").Value = 7051.0835
$ws.Range("K136").Value = 1906.1787
$ws.Range("L136").Value = 21153.2505
$ws.Range("M136").Value = 643.8212999999998
$ws.Range("N136").Value = -26253.2505
